$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.243.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.47%  "

$ws.Range("D3").Value = "'2.749.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.21%  "

$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").Value = "'353.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.18%  "

$ws.Range("D6").Value = "'106.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.54%  "

$ws.Range("E7").Value = "  -2.91%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").Value = "'0.577"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.72%  "

$ws.Range("D10").Value = "'39.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.47%  "

$ws.Range("E11").Value = "  +2.81%  "

$ws.Range("D12").Value = "'0.0833"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.05%  "

$ws.Range("D13").Value = "'19.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.72%  "

$ws.Range("E14").Value = "  -3.82%  "

$ws.Range("D15").Value = "'3.180.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.61%  "

$ws.Range("D16").Value = "'2.744.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.95%  "

$ws.Range("D17").Value = "'0.916"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").Value = "'51.153.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.69%  "

$ws.Range("D19").Value = "'7.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.92%  "

$ws.Range("D20").Value = "'3.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.73%  "

$ws.Range("D21").Value = "'12.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.49%  "

$ws.Range("D22").Value = "'0.0₃0958"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.35%  "

$ws.Range("D23").Value = "'69.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.05%  "

$ws.Range("D24").Value = "'264.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.43%  "

$ws.Range("D25").Value = "'2.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.97%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("D27").Value = "'25.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.71%  "

$ws.Range("E28").Value = "  +13.73%  "

$ws.Range("E29").Value = "  +0.73%  "

$ws.Range("D30").Value = "'10.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.91%  "

$ws.Range("E31").Value = "  +3.91%  "

$ws.Range("D32").Value = "'34.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.27%  "

$ws.Range("D33").Value = "'51.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.89%  "

$ws.Range("D34").Value = "'0.0440"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.78%  "

$ws.Range("D35").Value = "'0.0827"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.12%  "

$ws.Range("E36").Value = "  -6.24%  "

$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("D38").Value = "'18.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.40%  "

$ws.Range("D39").Value = "'3.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.50%  "

$ws.Range("E40").Value = "  -3.90%  "

$ws.Range("E41").Value = "  -3.00%  "

$ws.Range("D42").Value = "'2.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.68%  "

$ws.Range("D43").Value = "'2.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.46%  "

$ws.Range("D44").Value = "'119.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.27%  "

$ws.Range("D45").Value = "'21.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.52%  "

$ws.Range("D46").Value = "'2.070.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("D47").Value = "'2.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.44%  "

$ws.Range("E48").Value = "  -2.03%  "

$ws.Range("D49").Value = "'0.919"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.15%  "

$ws.Range("D50").Value = "'5.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.45%  "

$ws.Range("E51").Value = "  +5.05%  "
